$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Header block: "2020-03-10" -> "2020-02-04" (top date line, no
#    longer holds the _GoBack bookmark after this edit).
# ------------------------------------------------------------------
$d.Content.Find.Execute("2020-03-10", $true, $false, $false, $false, $false, `
    $true, 1, $false, "2020-02-04", 2) | Out-Null

# ------------------------------------------------------------------
# 2) "Last Updated: 2020-02-24" -> "Last Updated: 2020-03-09"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Last Updated: 2020-02-24", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Last Updated: 2020-03-09", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from the end of the date paragraph
#    down to the end of the "Last Updated: ..." paragraph.
# ------------------------------------------------------------------
$lastUpdatedPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Last Updated:*") {
        $lastUpdatedPara = $d.Paragraphs($i)
        break
    }
}
$pRange = $lastUpdatedPara.Range
$bmTarget = $d.Range($pRange.Start, $pRange.End - 1)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bmTarget)

# ------------------------------------------------------------------
# 4) "Raymond Knorr, Jaskirat Josan, Samuel Iregbu" paragraph: drop
#    the spell-check proofErr wrapping / run-splitting, keep a single
#    plain run with the same visible text.
# ------------------------------------------------------------------
function Set-ParagraphPlainText($paragraph, $text) {
    $r = $paragraph.Range
    $target = $d.Range($r.Start, $r.End - 1)
    $escaped = $text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($xml)
}

$namesPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Raymond Knorr*") {
        $namesPara = $d.Paragraphs($i)
        break
    }
}
Set-ParagraphPlainText $namesPara "Raymond Knorr, Jaskirat Josan, Samuel Iregbu"

# ------------------------------------------------------------------
# 5) Stakeholders paragraph: remove the gramStart/gramEnd proofErr
#    markers around "The" -- visible text is unchanged, but the tab
#    + text need to end up back-to-back in a single clean run.
# ------------------------------------------------------------------
$stakeholdersPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Stakeholders in the project would include*") {
        $stakeholdersPara = $d.Paragraphs($i)
        break
    }
}
$sRange = $stakeholdersPara.Range
$sTarget = $d.Range($sRange.Start, $sRange.End - 1)
$stakeholdersText = "Stakeholders in the project would include anyone that does/will have any vested interest in the project or outcome at any point in the development of the project. These people would include The JRS, Tim, Candace, Regina Fire and Protective Services (including all firefighters and workers), the U of R Software Engineering program. "
$sXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/><w:t xml:space="preserve">' + $stakeholdersText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$sTarget.InsertXML($sXml)
